$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a date-like string ("06/08/2025"); Excel's Value setter
# auto-parses date-looking text into a date serial. Force it to stay text
# by switching the cell to a text number format first, then clear the
# formatting again afterwards so the cell ends up unstyled, matching the
# other (unstyled) date cells in this column (A2, A3).
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "06/08/2025"
$ws.Range("A4").ClearFormats()

$ws.Range("B4").Value = "Basel"
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "Young Boys"
$ws.Range("F4").Value = "W"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2.81
$ws.Range("L4").Value = 0.47
$ws.Range("M4").Value = 22
$ws.Range("N4").Value = 13
$ws.Range("O4").Value = 12
$ws.Range("P4").Value = 4
